# Update the LR-pairs NATMI output sheet (Adam23-Itga4) with recomputed
# TPM-derived statistics (ligand/receptor expression, specificity and
# edge weight columns G..T) for data rows 2-10, per "update scripts wuth
# new tpm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value  = 0.07271233333333334
$ws.Range("H2").Value  = 0.218137
$ws.Range("I2").Value  = 0.004171225362010892
$ws.Range("J2").Value  = 0.004171225362010893
$ws.Range("K2").Value  = 2
$ws.Range("L2").Value  = 0.6666666666666666
$ws.Range("M2").Value  = 0.05968133333333333
$ws.Range("N2").Value  = 0.179044
$ws.Range("O2").Value  = 0.02602747651633847
$ws.Range("P2").Value  = 0.02602747651633848
$ws.Range("Q2").Value  = 0.004339569003111111
$ws.Range("R2").Value  = 0.039056121028
$ws.Range("S2").Value  = 0.000108566470154094
$ws.Range("T2").Value  = 0.000108566470154094

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value  = 0.07271233333333334
$ws.Range("H3").Value  = 0.218137
$ws.Range("I3").Value  = 0.004171225362010892
$ws.Range("J3").Value  = 0.004171225362010893
$ws.Range("O3").Value  = 0.144012433133819
$ws.Range("P3").Value  = 0.144012433133819
$ws.Range("Q3").Value  = 0.02401123637544444
$ws.Range("R3").Value  = 0.216101127379
$ws.Range("S3").Value  = 0.0006007083135326836
$ws.Range("T3").Value  = 0.0006007083135326838

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value  = 0.07271233333333334
$ws.Range("H4").Value  = 0.218137
$ws.Range("I4").Value  = 0.004171225362010892
$ws.Range("J4").Value  = 0.004171225362010893
$ws.Range("O4").Value  = 0.8299600903498424
$ws.Range("P4").Value  = 0.8299600903498425
$ws.Range("Q4").Value  = 0.1383794959776667
$ws.Range("R4").Value  = 1.245415463799
$ws.Range("S4").Value  = 0.003461950578324114
$ws.Range("T4").Value  = 0.003461950578324116

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value  = 0.5387060579248023
$ws.Range("J5").Value  = 0.5387060579248023
$ws.Range("K5").Value  = 2
$ws.Range("L5").Value  = 0.6666666666666666
$ws.Range("M5").Value  = 0.05968133333333333
$ws.Range("N5").Value  = 0.179044
$ws.Range("O5").Value  = 0.02602747651633847
$ws.Range("P5").Value  = 0.02602747651633848
$ws.Range("Q5").Value  = 0.5604473285115555
$ws.Range("R5").Value  = 5.044025956603999
$ws.Range("S5").Value  = 0.01402115927184707
$ws.Range("T5").Value  = 0.01402115927184707

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value  = 0.5387060579248023
$ws.Range("J6").Value  = 0.5387060579248023
$ws.Range("O6").Value  = 0.144012433133819
$ws.Range("P6").Value  = 0.144012433133819
$ws.Range("S6").Value  = 0.07758037014567883
$ws.Range("T6").Value  = 0.07758037014567883

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value  = 0.5387060579248023
$ws.Range("J7").Value  = 0.5387060579248023
$ws.Range("O7").Value  = 0.8299600903498424
$ws.Range("P7").Value  = 0.8299600903498425
$ws.Range("S7").Value  = 0.4471045285072764
$ws.Range("T7").Value  = 0.4471045285072764

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value  = 0.4571227167131868
$ws.Range("J8").Value  = 0.4571227167131868
$ws.Range("K8").Value  = 2
$ws.Range("L8").Value  = 0.6666666666666666
$ws.Range("M8").Value  = 0.05968133333333333
$ws.Range("N8").Value  = 0.179044
$ws.Range("O8").Value  = 0.02602747651633847
$ws.Range("P8").Value  = 0.02602747651633848
$ws.Range("Q8").Value  = 0.4755714208426667
$ws.Range("R8").Value  = 4.280142787583999
$ws.Range("S8").Value  = 0.01189775077433731
$ws.Range("T8").Value  = 0.01189775077433732

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value  = 0.4571227167131868
$ws.Range("J9").Value  = 0.4571227167131868
$ws.Range("O9").Value  = 0.144012433133819
$ws.Range("P9").Value  = 0.144012433133819
$ws.Range("Q9").Value  = 2.631380625834667
$ws.Range("S9").Value  = 0.0658313546746075
$ws.Range("T9").Value  = 0.0658313546746075

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.4571227167131868
$ws.Range("J10").Value = 0.4571227167131868
$ws.Range("O10").Value = 0.8299600903498424
$ws.Range("P10").Value = 0.8299600903498425
$ws.Range("S10").Value = 0.3793936112642419
$ws.Range("T10").Value = 0.379393611264242
